# Apply cryptos list update (prices & volume changes) as per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.268.41"
$ws.Range("E2").Value = "  -0.32%  "

$ws.Range("D3").Value = "3.438.73"
$ws.Range("E3").Value = "  +2.15%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "549.61"
$ws.Range("E5").Value = "  +2.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.75"
$ws.Range("E6").Value = "  -2.26%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.643"
$ws.Range("E7").Value = "  +6.72%  "

$ws.Range("E8").Value = "  -0.16%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.631"
$ws.Range("E9").Value = "  +0.75%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.150"
$ws.Range("E10").Value = "  +6.29%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.70"
$ws.Range("E11").Value = "  -2.64%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000271"
$ws.Range("E12").Value = "  +1.98%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.19"
$ws.Range("E13").Value = "  -0.20%  "

$ws.Range("D14").Value = "3.970.41"
$ws.Range("E14").Value = "  +1.71%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.121"
$ws.Range("E15").Value = "  +0.58%  "

$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.37"
$ws.Range("E16").Value = "  +2.27%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.426.59"
$ws.Range("E17").Value = "  +1.83%  "

$ws.Range("D18").Value = "65.159.83"
$ws.Range("E18").Value = "  -0.97%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.85"
$ws.Range("E19").Value = "  +3.59%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.985"
$ws.Range("E20").Value = "  +0.40%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "414.10"
$ws.Range("E21").Value = "  +6.02%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.00"
$ws.Range("E22").Value = "  +5.64%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "85.76"
$ws.Range("E23").Value = "  +3.29%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.27"
$ws.Range("E24").Value = "  +0.67%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.76"
$ws.Range("E25").Value = "  -8.91%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.86"
$ws.Range("E26").Value = "  +0.89%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.26"
$ws.Range("E27").Value = "  +6.18%  "

$ws.Range("E28").Value = "  +0.48%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.96"
$ws.Range("E29").Value = "  +6.09%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "29.81"
$ws.Range("E30").Value = "  +1.11%  "

$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.50"
$ws.Range("E31").Value = "  -4.00%  "

$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "607.71"
$ws.Range("E32").Value = "  -8.42%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.66"
$ws.Range("E33").Value = "  +2.05%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.109"
$ws.Range("E34").Value = "  +0.25%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "59.11"
$ws.Range("E35").Value = "  +1.96%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.11%  "

$ws.Range("E37").Value = "  +12.94%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "37.34"
$ws.Range("E38").Value = "  -0.98%  "

$ws.Range("D39").Value = "0.0₃0787"
$ws.Range("E39").Value = "  +1.08%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.378"
$ws.Range("E40").Value = "  -4.55%  "

$ws.Range("D41").Value = "3.194.09"
$ws.Range("E41").Value = "  +6.08%  "

$ws.Range("E42").Value = "  +0.30%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.998"
$ws.Range("E43").Value = "  -0.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.54"
$ws.Range("E44").Value = "  -8.24%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.79"
$ws.Range("E45").Value = "  +0.76%  "

$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0410"
$ws.Range("E46").Value = "  -0.46%  "

$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.20"
$ws.Range("E47").Value = "  +0.50%  "

$ws.Range("E48").Value = "  -0.79%  "

$ws.Range("E49").Value = "  +3.64%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "137.25"
$ws.Range("E50").Value = "  -1.28%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.36"
$ws.Range("E51").Value = "  -2.05%  "
